$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2262
$ws.Range("I70").Value = 2201
$ws.Range("J70").Value = 2272.1667
$ws.Range("K70").Value = 6603
$ws.Range("L70").Value = 6816.500100000001
$ws.Range("M70").Value = -6333
$ws.Range("N70").Value = -7356.500100000001
$ws.Range("H73").Value = 2262
$ws.Range("I73").Value = 2201
$ws.Range("J73").Value = 2272.1667
$ws.Range("K73").Value = 6603
$ws.Range("L73").Value = 6816.500100000001
$ws.Range("M73").Value = -5667
$ws.Range("N73").Value = -8688.500100000001
$ws.Range("H111").Value = 3021.0557
$ws.Range("I111").Value = 2631.5454
$ws.Range("K111").Value = 7894.6362
$ws.Range("M111").Value = -4827.6362
$ws.Range("H121").Value = 944.0714
$ws.Range("J121").Value = 976.9487
$ws.Range("L121").Value = 2930.8461
$ws.Range("N121").Value = -6424.846100000001
$ws.Range("H138").Value = 2673.17
$ws.Range("I138").Value = 1204.1111
$ws.Range("J138").Value = 3499.5156
$ws.Range("K138").Value = 3612.3333
$ws.Range("L138").Value = 10498.5468
$ws.Range("M138").Value = 1527.6667
$ws.Range("N138").Value = -20778.5468
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31476.73
$ws.Range("I32").Value = 11682.296
$ws.Range("J32").Value = 113699.766
$ws.Range("K32").Value = 11682.296
$ws.Range("L32").Value = 113699.766
$ws.Range("M32").Value = -11395.296
$ws.Range("N32").Value = -114273.766
$ws.Range("H61").Value = 2225.9546
$ws.Range("I61").Value = 2189.0952
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2189.0952
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1977.0952
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 5051.5
$ws.Range("I74").Value = 5182.9165
$ws.Range("J74").Value = 4893.8
$ws.Range("K74").Value = 5182.9165
$ws.Range("L74").Value = 4893.8
$ws.Range("M74").Value = -4308.9165
$ws.Range("N74").Value = -6641.8
$ws.Range("H77").Value = 5051.5
$ws.Range("I77").Value = 5182.9165
$ws.Range("J77").Value = 4893.8
$ws.Range("K77").Value = 25914.5825
$ws.Range("L77").Value = 24469
$ws.Range("M77").Value = -21546.5825
$ws.Range("N77").Value = -33205
$ws.Range("H136").Value = 2225.9546
$ws.Range("I136").Value = 2189.0952
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6567.285600000001
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4017.285600000001
$ws.Range("N136").Value = -14100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 10150
$ws.Range("J87").Value = 10150
$ws.Range("L87").Value = 10150
$ws.Range("N87").Value = -12646
$ws.Range("H90").Value = 10150
$ws.Range("J90").Value = 10150
$ws.Range("L90").Value = 30450
$ws.Range("N90").Value = -42930
$ws.Range("H124").Value = 42500
$ws.Range("J124").Value = 42500
$ws.Range("L124").Value = 42500
$ws.Range("N124").Value = -52320
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 21995.6
$ws.Range("J97").Value = 21995.6
$ws.Range("L97").Value = 21995.6
$ws.Range("N97").Value = -23977.6
$ws.Range("H132").Value = 2215.15
$ws.Range("I132").Value = 2156
$ws.Range("J132").Value = 2418.889
$ws.Range("K132").Value = 6468
$ws.Range("L132").Value = 7256.667
$ws.Range("M132").Value = -3938
$ws.Range("N132").Value = -12316.667
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1030.566
$ws.Range("I122").Value = 513.2222
$ws.Range("J122").Value = 1567.8077
$ws.Range("K122").Value = 4618.999800000001
$ws.Range("L122").Value = 14110.2693
$ws.Range("M122").Value = -2168.999800000001
$ws.Range("N122").Value = -19010.2693
$ws.Range("H134").Value = 3819.3125
$ws.Range("I134").Value = 3761
$ws.Range("J134").Value = 3916.5
$ws.Range("K134").Value = 11283
$ws.Range("L134").Value = 11749.5
$ws.Range("M134").Value = -6213
$ws.Range("N134").Value = -21889.5
$ws.Range("H140").Value = 26281.523
$ws.Range("J140").Value = 3369.1853
$ws.Range("L140").Value = 10107.5559
$ws.Range("N140").Value = -20467.5559
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H27").Value = 40000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H47").Value = 31677
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 31677
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 31677
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -32813
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 32500
$ws.Range("J63").Value = 32500
$ws.Range("L63").Value = 32500
$ws.Range("N63").Value = -33998
$ws.Range("H66").Value = 32500
$ws.Range("J66").Value = 32500
$ws.Range("L66").Value = 97500
$ws.Range("N66").Value = -104988
$ws.Range("H74").Value = 267497.5
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 353330
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 353330
$ws.Range("M74").Value = -9002
$ws.Range("N74").Value = -355326
$ws.Range("H77").Value = 267497.5
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 353330
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 1059990
$ws.Range("M77").Value = -25008
$ws.Range("N77").Value = -1069974
$ws.Range("H80").Value = 27641
$ws.Range("J80").Value = 27641
$ws.Range("L80").Value = 27641
$ws.Range("N80").Value = -29887
$ws.Range("H83").Value = 27641
$ws.Range("J83").Value = 27641
$ws.Range("L83").Value = 82923
$ws.Range("N83").Value = -94155
$ws.Range("H132").Value = 3923.175
$ws.Range("I132").Value = 4348.375
$ws.Range("J132").Value = 3285.375
$ws.Range("K132").Value = 13045.125
$ws.Range("L132").Value = 9856.125
$ws.Range("M132").Value = -10515.125
$ws.Range("N132").Value = -14916.125
$ws.Range("H136").Value = 5088.0645
$ws.Range("I136").Value = 1624.6
$ws.Range("J136").Value = 19519.166
$ws.Range("K136").Value = 4873.799999999999
$ws.Range("L136").Value = 58557.49800000001
$ws.Range("M136").Value = -2323.799999999999
$ws.Range("N136").Value = -63657.49800000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 28000
$ws.Range("J125").Value = 28000
$ws.Range("L125").Value = 28000
$ws.Range("N125").Value = -37840
$ws.Range("H132").Value = 3328.1562
$ws.Range("I132").Value = 4479.533
$ws.Range("J132").Value = 2312.2354
$ws.Range("K132").Value = 13438.599
$ws.Range("L132").Value = 6936.706200000001
$ws.Range("M132").Value = -10908.599
$ws.Range("N132").Value = -11996.7062
